$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number format / style from A301 down to the new date cells (A302:A328)
$ws.Range("A301").Copy() | Out-Null
$ws.Range("A302:A328").PasteSpecial(-4122) | Out-Null

$data = @(
  @(44376, 0, 5, 14.54926380725135),
  @(44377, 0, 5, 14.54926380725135),
  @(44378, 1, 5, 14.54926380725135),
  @(44379, 0, 5, 14.54926380725135),
  @(44380, 0, 2, 5.819705522900541),
  @(44381, 1, 3, 8.729558284350812),
  @(44382, 0, 2, 5.819705522900541),
  @(44383, 1, 3, 8.729558284350812),
  @(44384, 0, 3, 8.729558284350812),
  @(44385, 1, 3, 8.729558284350812),
  @(44386, 0, 3, 8.729558284350812),
  @(44387, 0, 3, 8.729558284350812),
  @(44388, 2, 4, 11.63941104580108),
  @(44389, 0, 4, 11.63941104580108),
  @(44390, 0, 3, 8.729558284350812),
  @(44391, 0, 3, 8.729558284350812),
  @(44392, 0, 2, 5.819705522900541),
  @(44393, 1, 3, 8.729558284350812),
  @(44394, 1, 4, 11.63941104580108),
  @(44395, 1, 3, 8.729558284350812),
  @(44396, 0, 3, 8.729558284350812),
  @(44397, 0, 3, 8.729558284350812),
  @(44398, 2, 5, 14.54926380725135),
  @(44399, 6, 11, 32.00838037595297),
  @(44400, 2, 12, 34.91823313740325),
  @(44401, 6, 17, 49.4674969446546),
  @(44402, 5, 21, 61.10690799045568)
)

$r = 302
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

Write-Output "done"
